$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a cell to be written as literal text (avoids Excel
# auto-converting numeric-looking strings like "1.000" or "0.9974"
# into real numbers), then strip the temporary formatting again so
# the cell keeps using the default (unstyled) cell format.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# --- Price / volume updates for unchanged rows ---
Set-TextValue $ws.Range("D2") "27.033.39"
$ws.Range("E2").Value = "  +0.31%  "
Set-TextValue $ws.Range("D3") "1.872.73"
$ws.Range("E3").Value = "  +0.63%  "
Set-TextValue $ws.Range("D4") "0.9974"
$ws.Range("E4").Value = "  -0.32%  "
Set-TextValue $ws.Range("D5") "305.64"
$ws.Range("E5").Value = "  -0.17%  "
Set-TextValue $ws.Range("D6") "0.9979"
$ws.Range("E6").Value = "  -0.28%  "
Set-TextValue $ws.Range("D7") "0.5056"
$ws.Range("E7").Value = "  +0.57%  "
Set-TextValue $ws.Range("D8") "0.3659"
$ws.Range("E8").Value = "  -1.56%  "
Set-TextValue $ws.Range("D9") "0.07184"
$ws.Range("E9").Value = "  +0.92%  "
Set-TextValue $ws.Range("D10") "0.8953"
$ws.Range("E10").Value = "  +1.33%  "
Set-TextValue $ws.Range("D11") "20.78"
$ws.Range("E11").Value = "  +1.30%  "
Set-TextValue $ws.Range("D14") "95.38"
$ws.Range("E14").Value = "  +7.14%  "
Set-TextValue $ws.Range("D15") "5.243"
$ws.Range("E15").Value = "  -0.78%  "
Set-TextValue $ws.Range("D16") "0.9974"
$ws.Range("E16").Value = "  -0.40%  "
Set-TextValue $ws.Range("D17") "0.000008542"
$ws.Range("E17").Value = "  +2.11%  "
Set-TextValue $ws.Range("D18") "14.25"
$ws.Range("E18").Value = "  +1.31%  "
Set-TextValue $ws.Range("D19") "1.0000"
$ws.Range("E19").Value = "  +0.01%  "
Set-TextValue $ws.Range("D20") "27.030.65"
$ws.Range("E20").Value = "  +0.14%  "
Set-TextValue $ws.Range("D21") "5.028"
$ws.Range("E21").Value = "  +0.10%  "
Set-TextValue $ws.Range("D22") "2.097.92"
$ws.Range("E22").Value = "  -0.70%  "
Set-TextValue $ws.Range("D23") "10.41"
$ws.Range("E23").Value = "  -0.37%  "
Set-TextValue $ws.Range("D24") "6.432"
$ws.Range("E24").Value = "  -0.34%  "
Set-TextValue $ws.Range("D25") "148.16"
$ws.Range("E25").Value = "  +0.72%  "
Set-TextValue $ws.Range("D26") "1.776"
$ws.Range("E26").Value = "  -3.92%  "
Set-TextValue $ws.Range("D27") "17.92"
$ws.Range("E27").Value = "  -0.05%  "
Set-TextValue $ws.Range("D28") "2.086"
$ws.Range("E28").Value = "  -0.51%  "
Set-TextValue $ws.Range("D29") "113.39"
$ws.Range("E29").Value = "  +0.76%  "
Set-TextValue $ws.Range("D30") "4.702"
$ws.Range("E30").Value = "  +1.40%  "
Set-TextValue $ws.Range("D31") "4.686"
$ws.Range("E31").Value = "  +0.55%  "
Set-TextValue $ws.Range("D32") "0.09180"
$ws.Range("E32").Value = "  +1.74%  "
Set-TextValue $ws.Range("D33") "0.05144"
$ws.Range("E33").Value = "  +0.52%  "
Set-TextValue $ws.Range("D34") "0.7552"
$ws.Range("E34").Value = "  +4.68%  "
Set-TextValue $ws.Range("D35") "2.971"
$ws.Range("E35").Value = "  -2.47%  "
Set-TextValue $ws.Range("D36") "1.158"
$ws.Range("E36").Value = "  +0.99%  "
Set-TextValue $ws.Range("D37") "3.204"
$ws.Range("E37").Value = "  +5.75%  "
Set-TextValue $ws.Range("D38") "2.589"
$ws.Range("E38").Value = "  +5.45%  "
Set-TextValue $ws.Range("D41") "1.072"
$ws.Range("E41").Value = "  +0.08%  "
Set-TextValue $ws.Range("D42") "6.591"
$ws.Range("E42").Value = "  +1.37%  "
Set-TextValue $ws.Range("D43") "116.27"
$ws.Range("E43").Value = "  +1.13%  "
Set-TextValue $ws.Range("D44") "8.593"
$ws.Range("E44").Value = "  +4.44%  "
Set-TextValue $ws.Range("D45") "0.1475"
$ws.Range("E45").Value = "  +0.90%  "
Set-TextValue $ws.Range("D46") "0.4728"
$ws.Range("E46").Value = "  +2.98%  "
Set-TextValue $ws.Range("D47") "0.9987"
$ws.Range("E47").Value = "  -0.20%  "
Set-TextValue $ws.Range("D48") "10.05"
$ws.Range("E48").Value = "  +0.90%  "
Set-TextValue $ws.Range("D49") "1.566"
$ws.Range("E49").Value = "  +0.44%  "
Set-TextValue $ws.Range("D50") "36.88"
$ws.Range("E50").Value = "  +1.24%  "
Set-TextValue $ws.Range("D51") "63.28"
$ws.Range("E51").Value = "  -0.74%  "

# --- Row 12/13 swap: TRON <-> WrappedEther (with updated prices) ---
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws.Range("D12") "1.861.51"
$ws.Range("E12").Value = "  +0.01%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D13") "0.07510"
$ws.Range("E13").Value = "  -0.74%  "

# --- Row 39/40 swap: VeChain <-> TheSandbox (with updated prices) ---
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D39") "0.5624"
$ws.Range("E39").Value = "  +6.71%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D40") "0.02002"
$ws.Range("E40").Value = "  -1.26%  "
